$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '66.946.44'
$ws.Range('E2').Value = '  +5.51%  '
$ws.Range('D3').Value = '3.704.21'
$ws.Range('E3').Value = '  +6.80%  '
Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  +0.07%  '
Set-TextValue $ws.Range('D5') '423.66'
$ws.Range('E5').Value = '  +2.03%  '
Set-TextValue $ws.Range('D6') '130.62'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('D7').Value = '3.696.11'
$ws.Range('E7').Value = '  +6.74%  '
Set-TextValue $ws.Range('D8') '0.642'
$ws.Range('E8').Value = '  +1.98%  '
$ws.Range('E9').Value = '  +0.11%  '
Set-TextValue $ws.Range('D10') '0.766'
$ws.Range('E10').Value = '  +1.69%  '
$ws.Range('E11').Value = '  +16.30%  '
Set-TextValue $ws.Range('D12') '0.0000387'
$ws.Range('E12').Value = '  +67.35%  '
Set-TextValue $ws.Range('D13') '42.76'
$ws.Range('E13').Value = '  +1.11%  '
Set-TextValue $ws.Range('D14') '10.07'
$ws.Range('E14').Value = '  +2.68%  '
$ws.Range('D15').Value = '4.294.53'
$ws.Range('E15').Value = '  +6.65%  '
$ws.Range('E16').Value = '  +0.06%  '
Set-TextValue $ws.Range('D17') '20.65'
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').Value = '3.729.69'
$ws.Range('E18').Value = '  +7.43%  '
Set-TextValue $ws.Range('D19') '12.90'
$ws.Range('E19').Value = '  +4.23%  '
$ws.Range('D21').Value = '67.000.92'
$ws.Range('E21').Value = '  +5.73%  '
Set-TextValue $ws.Range('D22') '446.14'
$ws.Range('E22').Value = '  -2.64%  '
Set-TextValue $ws.Range('D23') '15.55'
$ws.Range('E23').Value = '  +17.03%  '
Set-TextValue $ws.Range('D24') '89.11'
$ws.Range('E24').Value = '  -1.06%  '
$ws.Range('E25').Value = '  -4.36%  '
Set-TextValue $ws.Range('D26') '37.86'
$ws.Range('E26').Value = '  +13.11%  '
$ws.Range('B27').Value = 'Filecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D27') '10.27'
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D28') '3.30'
$ws.Range('E28').Value = '  +0.15%  '
Set-TextValue $ws.Range('D29') '5.00'
$ws.Range('E29').Value = '  +5.01%  '
$ws.Range('E30').Value = '  +1.51%  '
$ws.Range('E31').Value = '  +3.76%  '
$ws.Range('E32').Value = '  +6.92%  '
Set-TextValue $ws.Range('D33') '7.26'
$ws.Range('E33').Value = '  -3.39%  '
Set-TextValue $ws.Range('D34') '41.83'
$ws.Range('E34').Value = '  +4.48%  '
$ws.Range('E35').Value = '  -2.22%  '
$ws.Range('E36').Value = '  +0.04%  '
Set-TextValue $ws.Range('D37') '56.34'
$ws.Range('E37').Value = '  -2.62%  '
Set-TextValue $ws.Range('D38') '0.0490'
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('D39').Value = '0.0₃0733'
$ws.Range('E39').Value = '  +14.89%  '
Set-TextValue $ws.Range('D40') '2.99'
$ws.Range('E40').Value = '  +28.44%  '
$ws.Range('E41').Value = '  +6.03%  '
Set-TextValue $ws.Range('D42') '0.997'
$ws.Range('E42').Value = '  -0.24%  '
Set-TextValue $ws.Range('D43') '3.45'
$ws.Range('E43').Value = '  +3.29%  '
Set-TextValue $ws.Range('D44') '27.67'
$ws.Range('E44').Value = '  +26.14%  '
$ws.Range('E45').Value = '  +5.64%  '
$ws.Range('E46').Value = '  -5.05%  '
Set-TextValue $ws.Range('D47') '145.25'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D48') '2.64'
$ws.Range('E48').Value = '  -7.09%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D49') '4.34'
$ws.Range('E49').Value = '  -2.35%  '
Set-TextValue $ws.Range('D50') '0.306'
$ws.Range('E50').Value = '  -3.74%  '
Set-TextValue $ws.Range('D51') '0.158'
$ws.Range('E51').Value = '  +13.83%  '
